$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 351.34784
$ws.Range("I33").Value = 289.7143
$ws.Range("K33").Value = 289.7143
$ws.Range("M33").Value = -60.71429999999998
$ws.Range("H137").Value = 10167.694
$ws.Range("I137").Value = 2902.4546
$ws.Range("J137").Value = 13364.4
$ws.Range("K137").Value = 8707.363799999999
$ws.Range("L137").Value = 40093.2
$ws.Range("M137").Value = -6157.363799999999
$ws.Range("N137").Value = -45193.2
$ws.Range("H138").Value = 6088.76
$ws.Range("J138").Value = 5847
$ws.Range("L138").Value = 17541
$ws.Range("N138").Value = -27821

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1513764.9
$ws.Range("I32").Value = 687644.1
$ws.Range("K32").Value = 687644.1
$ws.Range("M32").Value = -687357.1
$ws.Range("H37").Value = 40829.6
$ws.Range("I37").Value = 25017
$ws.Range("J37").Value = 51371.332
$ws.Range("K37").Value = 25017
$ws.Range("L37").Value = 51371.332
$ws.Range("M37").Value = -24744
$ws.Range("N37").Value = -51917.332
$ws.Range("H44").Value = 77000
$ws.Range("J44").Value = 77000
$ws.Range("L44").Value = 77000
$ws.Range("N44").Value = -77976
$ws.Range("H55").Value = 59009.668
$ws.Range("J55").Value = 77000
$ws.Range("L55").Value = 77000
$ws.Range("N55").Value = -77630
$ws.Range("H74").Value = 23811178
$ws.Range("I74").Value = 1322.9259
$ws.Range("J74").Value = 66668916
$ws.Range("K74").Value = 1322.9259
$ws.Range("L74").Value = 66668916
$ws.Range("M74").Value = -448.9259
$ws.Range("N74").Value = -66670664
$ws.Range("H77").Value = 23811178
$ws.Range("I77").Value = 1322.9259
$ws.Range("J77").Value = 66668916
$ws.Range("K77").Value = 6614.6295
$ws.Range("L77").Value = 333344580
$ws.Range("M77").Value = -2246.6295
$ws.Range("N77").Value = -333353316
$ws.Range("H80").Value = 88633.336
$ws.Range("I80").Value = 78000
$ws.Range("J80").Value = 93950
$ws.Range("K80").Value = 78000
$ws.Range("L80").Value = 93950
$ws.Range("M80").Value = -77002
$ws.Range("N80").Value = -95946
$ws.Range("H83").Value = 88633.336
$ws.Range("I83").Value = 78000
$ws.Range("J83").Value = 93950
$ws.Range("K83").Value = 234000
$ws.Range("L83").Value = 281850
$ws.Range("M83").Value = -229008
$ws.Range("N83").Value = -291834
$ws.Range("H88").Value = 4434.5
$ws.Range("I88").Value = 2360.4
$ws.Range("J88").Value = 5916
$ws.Range("K88").Value = 2360.4
$ws.Range("L88").Value = 5916
$ws.Range("M88").Value = -1954.4
$ws.Range("N88").Value = -6728
$ws.Range("H91").Value = 4434.5
$ws.Range("I91").Value = 2360.4
$ws.Range("J91").Value = 5916
$ws.Range("K91").Value = 2360.4
$ws.Range("L91").Value = 5916
$ws.Range("M91").Value = -956.4000000000001
$ws.Range("N91").Value = -8724

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 627
$ws.Range("I134").Value = 614.6429000000001
$ws.Range("J134").Value = 800
$ws.Range("K134").Value = 1843.9287
$ws.Range("L134").Value = 2400
$ws.Range("M134").Value = 691.0712999999998
$ws.Range("N134").Value = -7470

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2196507.5
$ws.Range("I31").Value = 3287.4443
$ws.Range("J31").Value = 2607736.2
$ws.Range("K31").Value = 3287.4443
$ws.Range("L31").Value = 2607736.2
$ws.Range("M31").Value = -2992.4443
$ws.Range("N31").Value = -2608326.2
$ws.Range("H34").Value = 2196507.5
$ws.Range("I34").Value = 3287.4443
$ws.Range("J34").Value = 2607736.2
$ws.Range("K34").Value = 3287.4443
$ws.Range("L34").Value = 2607736.2
$ws.Range("M34").Value = -3085.4443
$ws.Range("N34").Value = -2608140.2
$ws.Range("H58").Value = 3251.375
$ws.Range("I58").Value = 2252.75
$ws.Range("K58").Value = 2252.75
$ws.Range("M58").Value = -2049.75
$ws.Range("H99").Value = 142862060
$ws.Range("I99").Value = 250003380
$ws.Range("J99").Value = 6999.6665
$ws.Range("K99").Value = 250003380
$ws.Range("L99").Value = 6999.6665
$ws.Range("M99").Value = -250001882
$ws.Range("N99").Value = -9995.666499999999
$ws.Range("H126").Value = 142862060
$ws.Range("I126").Value = 250003380
$ws.Range("J126").Value = 6999.6665
$ws.Range("K126").Value = 750010140
$ws.Range("L126").Value = 20998.9995
$ws.Range("M126").Value = -750007670
$ws.Range("N126").Value = -25938.9995
$ws.Range("H136").Value = 3251.375
$ws.Range("I136").Value = 2252.75
$ws.Range("K136").Value = 6758.25
$ws.Range("M136").Value = -4208.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1337599.1
$ws.Range("I68").Value = 2874.3845
$ws.Range("J68").Value = 2045820.4
$ws.Range("K68").Value = 8623.1535
$ws.Range("L68").Value = 6137461.199999999
$ws.Range("M68").Value = -7812.1535
$ws.Range("N68").Value = -6139083.199999999
$ws.Range("H71").Value = 1337599.1
$ws.Range("I71").Value = 2874.3845
$ws.Range("J71").Value = 2045820.4
$ws.Range("K71").Value = 25869.4605
$ws.Range("L71").Value = 18412383.6
$ws.Range("M71").Value = -21813.4605
$ws.Range("N71").Value = -18420495.6
$ws.Range("H107").Value = 4857.5713
$ws.Range("I107").Value = 999
$ws.Range("J107").Value = 10002.333
$ws.Range("K107").Value = 2997
$ws.Range("L107").Value = 30006.999
$ws.Range("M107").Value = -1077
$ws.Range("N107").Value = -33846.999
$ws.Range("H113").Value = 755.5599999999999
$ws.Range("I113").Value = 745
$ws.Range("J113").Value = 758.2
$ws.Range("K113").Value = 2235
$ws.Range("L113").Value = 2274.6
$ws.Range("M113").Value = -65
$ws.Range("N113").Value = -6614.6
$ws.Range("H131").Value = 2274954.5
$ws.Range("I131").Value = 19277.125
$ws.Range("K131").Value = 57831.375
$ws.Range("M131").Value = -52791.375
$ws.Range("H137").Value = 3632.818
$ws.Range("J137").Value = 3945.3333
$ws.Range("L137").Value = 11835.9999
$ws.Range("N137").Value = -22035.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2416.2
$ws.Range("I97").Value = 1860.3334
$ws.Range("J97").Value = 3250
$ws.Range("K97").Value = 1860.3334
$ws.Range("L97").Value = 3250
$ws.Range("M97").Value = -1364.3334
$ws.Range("N97").Value = -4242

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1616.6666
$ws.Range("I22").Value = 2233.3333
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2233.3333
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1938.3333
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1616.6666
$ws.Range("I27").Value = 2233.3333
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 2233.3333
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -2126.3333
$ws.Range("N27").Value = -1214
$ws.Range("H63").Value = 49979
$ws.Range("J63").Value = 49969.5
$ws.Range("L63").Value = 49969.5
$ws.Range("N63").Value = -51467.5
$ws.Range("H66").Value = 49979
$ws.Range("J66").Value = 49969.5
$ws.Range("L66").Value = 149908.5
$ws.Range("N66").Value = -157396.5
$ws.Range("H122").Value = 3485.7778
$ws.Range("I122").Value = 3821.8572
$ws.Range("J122").Value = 2309.5
$ws.Range("K122").Value = 11465.5716
$ws.Range("L122").Value = 6928.5
$ws.Range("M122").Value = -9015.571599999999
$ws.Range("N122").Value = -11828.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 32511.334
$ws.Range("J54").Value = 43500
$ws.Range("L54").Value = 43500
$ws.Range("N54").Value = -44540
$ws.Range("H132").Value = 1651.4667
$ws.Range("I132").Value = 1690.1538
$ws.Range("K132").Value = 5070.4614
$ws.Range("M132").Value = -2540.4614
